# Refresh market-price-derived columns (H-N) on several Leve sheets.
# Values below come from the scheduled market-data runner update.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H76").Value = 3465.5833
$ws.Range("I76").Value = 3399.889
$ws.Range("K76").Value = 3399.889
$ws.Range("M76").Value = -3084.889
$ws.Range("H79").Value = 3465.5833
$ws.Range("I79").Value = 3399.889
$ws.Range("K79").Value = 3399.889
$ws.Range("M79").Value = -2307.889
$ws.Range("H98").Value = 1125.8823
$ws.Range("I98").Value = 1254.6666
$ws.Range("J98").Value = 629.1429000000001
$ws.Range("K98").Value = 1254.6666
$ws.Range("L98").Value = 629.1429000000001
$ws.Range("M98").Value = 243.3334
$ws.Range("N98").Value = -3625.1429
$ws.Range("H113").Value = 2704.1482
$ws.Range("I113").Value = 2640
$ws.Range("J113").Value = 2887.4285
$ws.Range("K113").Value = 2640
$ws.Range("L113").Value = 2887.4285
$ws.Range("M113").Value = 614
$ws.Range("N113").Value = -9395.4285
$ws.Range("H116").Value = 1936.3125
$ws.Range("I116").Value = 1700
$ws.Range("J116").Value = 2015.0834
$ws.Range("K116").Value = 1700
$ws.Range("L116").Value = 2015.0834
$ws.Range("M116").Value = 1742
$ws.Range("N116").Value = -8899.0834
$ws.Range("H122").Value = 1125.8823
$ws.Range("I122").Value = 1254.6666
$ws.Range("J122").Value = 629.1429000000001
$ws.Range("K122").Value = 3763.9998
$ws.Range("L122").Value = 1887.4287
$ws.Range("M122").Value = -1313.9998
$ws.Range("N122").Value = -6787.4287
$ws.Range("H127").Value = 1394.1052
$ws.Range("I127").Value = 549.25
$ws.Range("J127").Value = 2008.5454
$ws.Range("K127").Value = 1647.75
$ws.Range("L127").Value = 6025.6362
$ws.Range("M127").Value = 3312.25
$ws.Range("N127").Value = -15945.6362
$ws.Range("H129").Value = 992.48
$ws.Range("J129").Value = 1120.6666
$ws.Range("L129").Value = 3361.9998
$ws.Range("N129").Value = -13361.9998
$ws.Range("H132").Value = 1326981.6
$ws.Range("I132").Value = 2576.1667
$ws.Range("J132").Value = 7003005
$ws.Range("K132").Value = 7728.500100000001
$ws.Range("L132").Value = 21009015
$ws.Range("M132").Value = -5198.500100000001
$ws.Range("N132").Value = -21014075

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2843.5715
$ws.Range("I2").Value = 3087.2856
$ws.Range("K2").Value = 3087.2856
$ws.Range("M2").Value = -2974.2856
$ws.Range("H61").Value = 91091860
$ws.Range("I61").Value = 111222936
$ws.Range("J61").Value = 502000
$ws.Range("K61").Value = 111222936
$ws.Range("L61").Value = 502000
$ws.Range("M61").Value = -111222724
$ws.Range("N61").Value = -502424
$ws.Range("H74").Value = 7201629.5
$ws.Range("I74").Value = 10041385
$ws.Range("K74").Value = 10041385
$ws.Range("M74").Value = -10040511
$ws.Range("H77").Value = 7201629.5
$ws.Range("I77").Value = 10041385
$ws.Range("K77").Value = 50206925
$ws.Range("M77").Value = -50202557
$ws.Range("H116").Value = 2843.5715
$ws.Range("I116").Value = 3087.2856
$ws.Range("K116").Value = 3087.2856
$ws.Range("M116").Value = -793.2856000000002
$ws.Range("H122").Value = 4117385.2
$ws.Range("I122").Value = 2245.0908
$ws.Range("J122").Value = 22224002
$ws.Range("K122").Value = 6735.2724
$ws.Range("L122").Value = 66672006
$ws.Range("M122").Value = -4285.2724
$ws.Range("N122").Value = -66676906
$ws.Range("H132").Value = 70749.53
$ws.Range("I132").Value = 44455.914
$ws.Range("J132").Value = 157142.86
$ws.Range("K132").Value = 133367.742
$ws.Range("L132").Value = 471428.58
$ws.Range("M132").Value = -130837.742
$ws.Range("N132").Value = -476488.58
$ws.Range("H136").Value = 91091860
$ws.Range("I136").Value = 111222936
$ws.Range("J136").Value = 502000
$ws.Range("K136").Value = 333668808
$ws.Range("L136").Value = 1506000
$ws.Range("M136").Value = -333666258
$ws.Range("N136").Value = -1511100

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2843.5715
$ws.Range("I3").Value = 3087.2856
$ws.Range("K3").Value = 3087.2856
$ws.Range("M3").Value = -2973.2856
$ws.Range("H99").Value = 1138.909
$ws.Range("I99").Value = 1003.1111
$ws.Range("K99").Value = 1003.1111
$ws.Range("M99").Value = 494.8889
$ws.Range("H134").Value = 2351.2222
$ws.Range("I134").Value = 1438.9375
$ws.Range("J134").Value = 3081.05
$ws.Range("K134").Value = 4316.8125
$ws.Range("L134").Value = 9243.150000000001
$ws.Range("M134").Value = -1781.8125
$ws.Range("N134").Value = -14313.15

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2551.5
$ws.Range("I62").Value = 2066.6667
$ws.Range("K62").Value = 2066.6667
$ws.Range("M62").Value = -1442.6667
$ws.Range("H65").Value = 2551.5
$ws.Range("I65").Value = 2066.6667
$ws.Range("K65").Value = 10333.3335
$ws.Range("M65").Value = -7213.333500000001
$ws.Range("H99").Value = 6254.793
$ws.Range("I99").Value = 6937.5
$ws.Range("J99").Value = 5137.636
$ws.Range("K99").Value = 6937.5
$ws.Range("L99").Value = 5137.636
$ws.Range("M99").Value = -5439.5
$ws.Range("N99").Value = -8133.636
$ws.Range("H126").Value = 6254.793
$ws.Range("I126").Value = 6937.5
$ws.Range("J126").Value = 5137.636
$ws.Range("K126").Value = 20812.5
$ws.Range("L126").Value = 15412.908
$ws.Range("M126").Value = -18342.5
$ws.Range("N126").Value = -20352.908
$ws.Range("H132").Value = 58783.37
$ws.Range("I132").Value = 39987.617
$ws.Range("J132").Value = 113082.22
$ws.Range("K132").Value = 119962.851
$ws.Range("L132").Value = 339246.66
$ws.Range("M132").Value = -117432.851
$ws.Range("N132").Value = -344306.66
$ws.Range("H134").Value = 26103.363
$ws.Range("I134").Value = 1672.8387
$ws.Range("J134").Value = 84360.766
$ws.Range("K134").Value = 5018.5161
$ws.Range("L134").Value = 253082.298
$ws.Range("M134").Value = -2483.5161
$ws.Range("N134").Value = -258152.298

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -2831
$ws.Range("N17").Value = -3338
$ws.Range("H34").Value = 1095.0769
$ws.Range("J34").Value = 1467
$ws.Range("L34").Value = 4401
$ws.Range("N34").Value = -4569
$ws.Range("H39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 2051
$ws.Range("J55").Value = 3800
$ws.Range("L55").Value = 11400
$ws.Range("N55").Value = -11754
$ws.Range("H113").Value = 475.5102
$ws.Range("I113").Value = 440.625
$ws.Range("J113").Value = 541.17645
$ws.Range("K113").Value = 1321.875
$ws.Range("L113").Value = 1623.52935
$ws.Range("M113").Value = 848.125
$ws.Range("N113").Value = -5963.529350000001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 942.8
$ws.Range("I102").Value = 925
$ws.Range("J102").Value = 1014
$ws.Range("K102").Value = 925
$ws.Range("L102").Value = 1014
$ws.Range("M102").Value = 697
$ws.Range("N102").Value = -4258
$ws.Range("H113").Value = 1944.4445
$ws.Range("J113").Value = 2037.5
$ws.Range("L113").Value = 2037.5
$ws.Range("N113").Value = -6377.5
$ws.Range("H132").Value = 75677.03999999999
$ws.Range("I132").Value = 57231.5
$ws.Range("J132").Value = 112568.11
$ws.Range("K132").Value = 171694.5
$ws.Range("L132").Value = 337704.33
$ws.Range("M132").Value = -169164.5
$ws.Range("N132").Value = -342764.33

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3680.8
$ws.Range("I61").Value = 3401.3333
$ws.Range("J61").Value = 4100
$ws.Range("K61").Value = 3401.3333
$ws.Range("L61").Value = 4100
$ws.Range("M61").Value = -3199.3333
$ws.Range("N61").Value = -4504
$ws.Range("H93").Value = 500
$ws.Range("I93").Value = 500
$ws.Range("K93").Value = 500
$ws.Range("M93").Value = 748
$ws.Range("H113").Value = 3680.8
$ws.Range("I113").Value = 3401.3333
$ws.Range("J113").Value = 4100
$ws.Range("K113").Value = 3401.3333
$ws.Range("L113").Value = 4100
$ws.Range("M113").Value = -1231.3333
$ws.Range("N113").Value = -8440
